$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs" (C) and "balls" (D) columns are stored as text in this sheet,
# so prefix the new values with an apostrophe to keep Excel from coercing
# them to numbers, then restore the default "Normal" style so the
# quote-prefix indicator doesn't leave a stray per-cell format behind.

$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'0"

$ws.Range("C3").Value = "'3"
$ws.Range("D3").Value = "'2"

$ws.Range("C4").Value = "'3"
$ws.Range("D4").Value = "'1"

$ws.Range("C5").Value = "'1"
$ws.Range("D5").Value = "'3"

$ws.Range("C6").Value = "'0"
$ws.Range("D6").Value = "'0"

$ws.Range("C2:D6").Style = "Normal"
